$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$ws.Cells.Item(104, 1).Value = 'Demo'
$ws.Cells.Item(104, 2).Value = 42577
$ws.Cells.Item(104, 3).Value = '1600'
$ws.Cells.Item(104, 4).Value = 'SSB'
$ws.Cells.Item(104, 5).Value = 'S125'

$ws.Cells.Item(105, 1).Value = 'Demo'
$ws.Cells.Item(105, 2).Value = 42577
$ws.Cells.Item(105, 3).Value = '1730'
$ws.Cells.Item(105, 4).Value = 'SSB'
$ws.Cells.Item(105, 5).Value = 'S125'

$ws.Cells.Item(106, 1).Value = 'Crestron Logout'
$ws.Cells.Item(106, 2).Value = 42577
$ws.Cells.Item(106, 3).Value = '1630'
$ws.Cells.Item(106, 4).Value = 'OSG'
$ws.Cells.Item(106, 5).Value = '1001'

$ws.Cells.Item(107, 1).Value = 'Crestron Logout'
$ws.Cells.Item(107, 2).Value = 42577
$ws.Cells.Item(107, 3).Value = '1630'
$ws.Cells.Item(107, 4).Value = 'OSG'
$ws.Cells.Item(107, 5).Value = '1003'

$ws.Cells.Item(108, 1).Value = 'Crestron Logout'
$ws.Cells.Item(108, 2).Value = 42577
$ws.Cells.Item(108, 3).Value = '1630'
$ws.Cells.Item(108, 4).Value = 'OSG'
$ws.Cells.Item(108, 5).Value = '2002'

$ws.Cells.Item(109, 1).Value = 'Crestron Logout'
$ws.Cells.Item(109, 2).Value = 42577
$ws.Cells.Item(109, 3).Value = '1600'
$ws.Cells.Item(109, 4).Value = 'OSG'
$ws.Cells.Item(109, 5).Value = '2027'

$ws.Cells.Item(110, 1).Value = 'Crestron Logout'
$ws.Cells.Item(110, 2).Value = 42577
$ws.Cells.Item(110, 3).Value = '1800'
$ws.Cells.Item(110, 4).Value = 'SSB'
$ws.Cells.Item(110, 5).Value = 'S125'

$ws.Cells.Item(111, 1).Value = 'Crestron Logout'
$ws.Cells.Item(111, 2).Value = 42577
$ws.Cells.Item(111, 3).Value = '2030'
$ws.Cells.Item(111, 4).Value = 'SSB'
$ws.Cells.Item(111, 5).Value = 'S123'

$ws.Cells.Item(112, 1).Value = 'Crestron Logout'
$ws.Cells.Item(112, 2).Value = 42577
$ws.Cells.Item(112, 3).Value = '1900'
$ws.Cells.Item(112, 4).Value = 'ACE'
$ws.Cells.Item(112, 5).Value = '002'

$ws.Cells.Item(113, 1).Value = 'Crestron Logout'
$ws.Cells.Item(113, 2).Value = 42577
$ws.Cells.Item(113, 3).Value = '1900'
$ws.Cells.Item(113, 4).Value = 'ACE'
$ws.Cells.Item(113, 5).Value = '003'

$ws.Cells.Item(114, 1).Value = 'Crestron Logout'
$ws.Cells.Item(114, 2).Value = 42577
$ws.Cells.Item(114, 3).Value = '1730'
$ws.Cells.Item(114, 4).Value = 'ACE'
$ws.Cells.Item(114, 5).Value = '006'

$ws.Cells.Item(115, 1).Value = 'Crestron Logout'
$ws.Cells.Item(115, 2).Value = 42577
$ws.Cells.Item(115, 3).Value = '1900'
$ws.Cells.Item(115, 4).Value = 'ACE'
$ws.Cells.Item(115, 5).Value = '009'

$ws.Cells.Item(116, 1).Value = 'Crestron Logout'
$ws.Cells.Item(116, 2).Value = 42577
$ws.Cells.Item(116, 3).Value = '1900'
$ws.Cells.Item(116, 4).Value = 'ACE'
$ws.Cells.Item(116, 5).Value = '010'

$ws.Cells.Item(117, 1).Value = 'Crestron Logout'
$ws.Cells.Item(117, 2).Value = 42577
$ws.Cells.Item(117, 3).Value = '1900'
$ws.Cells.Item(117, 4).Value = 'ACE'
$ws.Cells.Item(117, 5).Value = '011'

$ws.Cells.Item(118, 1).Value = 'Crestron Logout'
$ws.Cells.Item(118, 2).Value = 42577
$ws.Cells.Item(118, 3).Value = '2000'
$ws.Cells.Item(118, 4).Value = 'ACE'
$ws.Cells.Item(118, 5).Value = '013'

$ws.Cells.Item(119, 1).Value = 'Crestron Logout'
$ws.Cells.Item(119, 2).Value = 42577
$ws.Cells.Item(119, 3).Value = '1700'
$ws.Cells.Item(119, 4).Value = 'OSG'
$ws.Cells.Item(119, 5).Value = '1002'

$ws.Cells.Item(120, 1).Value = 'Crestron Logout'
$ws.Cells.Item(120, 2).Value = 42577
$ws.Cells.Item(120, 3).Value = '1700'
$ws.Cells.Item(120, 4).Value = 'OSG'
$ws.Cells.Item(120, 5).Value = '2001'

$ws.Cells.Item(121, 1).Value = 'Crestron Logout'
$ws.Cells.Item(121, 2).Value = 42577
$ws.Cells.Item(121, 3).Value = '1700'
$ws.Cells.Item(121, 4).Value = 'OSG'
$ws.Cells.Item(121, 5).Value = '2010'

$ws.Cells.Item(122, 1).Value = 'Crestron Logout'
$ws.Cells.Item(122, 2).Value = 42577
$ws.Cells.Item(122, 3).Value = '1730'
$ws.Cells.Item(122, 4).Value = 'TEL'
$ws.Cells.Item(122, 5).Value = '0001'

$ws.Cells.Item(123, 1).Value = 'Crestron Logout'
$ws.Cells.Item(123, 2).Value = 42577
$ws.Cells.Item(123, 3).Value = '1730'
$ws.Cells.Item(123, 4).Value = 'TEL'
$ws.Cells.Item(123, 5).Value = '0004'

$ws.Cells.Item(124, 1).Value = 'Crestron Logout'
$ws.Cells.Item(124, 2).Value = 42577
$ws.Cells.Item(124, 3).Value = '1700'
$ws.Cells.Item(124, 4).Value = 'SSB'
$ws.Cells.Item(124, 5).Value = 'N105'

$ws.Cells.Item(125, 1).Value = 'Crestron Logout'
$ws.Cells.Item(125, 2).Value = 42577
$ws.Cells.Item(125, 3).Value = '1700'
$ws.Cells.Item(125, 4).Value = 'SSB'
$ws.Cells.Item(125, 5).Value = 'N106'

$ws.Cells.Item(126, 1).Value = 'Crestron Logout'
$ws.Cells.Item(126, 2).Value = 42577
$ws.Cells.Item(126, 3).Value = '2030'
$ws.Cells.Item(126, 4).Value = 'SSB'
$ws.Cells.Item(126, 5).Value = 'N107'

$ws.Cells.Item(127, 1).Value = 'Crestron Logout'
$ws.Cells.Item(127, 2).Value = 42577
$ws.Cells.Item(127, 3).Value = '2030'
$ws.Cells.Item(127, 4).Value = 'SSB'
$ws.Cells.Item(127, 5).Value = 'N108'

$ws.Cells.Item(128, 1).Value = 'Crestron Logout'
$ws.Cells.Item(128, 2).Value = 42577
$ws.Cells.Item(128, 3).Value = '1700'
$ws.Cells.Item(128, 4).Value = 'SSB'
$ws.Cells.Item(128, 5).Value = 'N109'

$ws.Cells.Item(129, 1).Value = 'Crestron Logout'
$ws.Cells.Item(129, 2).Value = 42577
$ws.Cells.Item(129, 3).Value = '1600'
$ws.Cells.Item(129, 4).Value = 'SSB'
$ws.Cells.Item(129, 5).Value = 'N201'

$ws.Cells.Item(130, 1).Value = 'Crestron Logout'
$ws.Cells.Item(130, 2).Value = 42577
$ws.Cells.Item(130, 3).Value = '1900'
$ws.Cells.Item(130, 4).Value = 'SSB'
$ws.Cells.Item(130, 5).Value = 'S127'

$ws.Cells.Item(131, 1).Value = 'Crestron Logout'
$ws.Cells.Item(131, 2).Value = 42577
$ws.Cells.Item(131, 3).Value = '1900'
$ws.Cells.Item(131, 4).Value = 'SSB'
$ws.Cells.Item(131, 5).Value = 'S129'

$ws.Cells.Item(132, 1).Value = 'Crestron Logout'
$ws.Cells.Item(132, 2).Value = 42577
$ws.Cells.Item(132, 3).Value = '1900'
$ws.Cells.Item(132, 4).Value = 'SSB'
$ws.Cells.Item(132, 5).Value = 'W132'

$ws.Cells.Item(133, 1).Value = 'Crestron Logout'
$ws.Cells.Item(133, 2).Value = 42577
$ws.Cells.Item(133, 3).Value = '1900'
$ws.Cells.Item(133, 4).Value = 'SSB'
$ws.Cells.Item(133, 5).Value = 'W255'

$ws.Cells.Item(134, 1).Value = 'Crestron Logout'
$ws.Cells.Item(134, 2).Value = 42577
$ws.Cells.Item(134, 3).Value = '1630'
$ws.Cells.Item(134, 4).Value = 'SSB'
$ws.Cells.Item(134, 5).Value = 'W356'

$ws.Range("C134").Select() | Out-Null
